$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country")

# Populate column B with Natural Earth country name mapping
$ws.Cells.Item(1, 2).Value = "Natural Earth"
$ws.Cells.Item(2, 2).Value = "Albania"
$ws.Cells.Item(3, 2).Value = "Algeria"
$ws.Cells.Item(4, 2).Value = "Andorra"
$ws.Cells.Item(5, 2).Value = "Armenia"
$ws.Cells.Item(6, 2).Value = "Austria"
$ws.Cells.Item(7, 2).Value = "Azerbaijan"
$ws.Cells.Item(8, 2).Value = "Bahrain"
$ws.Cells.Item(9, 2).Value = "Belarus"
$ws.Cells.Item(10, 2).Value = "Belgium"
$ws.Cells.Item(11, 2).Value = "Bosnia and Herzegovina"
$ws.Cells.Item(12, 2).Value = "Bulgaria"
$ws.Cells.Item(13, 2).Value = "Cape Verde"
$ws.Cells.Item(14, 2).Value = "Croatia"
$ws.Cells.Item(15, 2).Value = "Cyprus"
$ws.Cells.Item(16, 2).Value = "Czech Republic"
$ws.Cells.Item(17, 2).Value = "Denmark"
$ws.Cells.Item(18, 2).Value = "Egypt"
$ws.Cells.Item(19, 2).Value = "Eritrea"
$ws.Cells.Item(20, 2).Value = "Estonia"
$ws.Cells.Item(21, 2).Value = "Ethiopia"
$ws.Cells.Item(22, 2).Value = "Faroe Islands"
$ws.Cells.Item(23, 2).Value = "Finland"
$ws.Cells.Item(24, 2).Value = "France"
$ws.Cells.Item(25, 2).Value = "Georgia"
$ws.Cells.Item(26, 2).Value = "Germany"
$ws.Cells.Item(27, 2).Value = "Gibraltar"
$ws.Cells.Item(28, 2).Value = "Greece"
$ws.Cells.Item(29, 2).Value = "Greenland"
$ws.Cells.Item(30, 2).Value = "Guernsey"
$ws.Cells.Item(31, 2).Value = "Hungary"
$ws.Cells.Item(32, 2).Value = "Iceland"
$ws.Cells.Item(33, 2).Value = "Iran"
$ws.Cells.Item(34, 2).Value = "Iraq"
$ws.Cells.Item(35, 2).Value = "Ireland"
$ws.Cells.Item(36, 2).Value = "Isle of Man"
$ws.Cells.Item(37, 2).Value = "Israel"
$ws.Cells.Item(38, 2).Value = "Italy"
$ws.Cells.Item(39, 2).Value = "India"
$ws.Cells.Item(40, 2).Value = "Japan"
$ws.Cells.Item(41, 2).Value = "Jersey"
$ws.Cells.Item(42, 2).Value = "Kazakhstan"
$ws.Cells.Item(43, 2).Value = "Kuwait"
$ws.Cells.Item(44, 2).Value = "Kyrgyzstan"
$ws.Cells.Item(45, 2).Value = "Latvia"
$ws.Cells.Item(46, 2).Value = "Lebanon"
$ws.Cells.Item(47, 2).Value = "Libya"
$ws.Cells.Item(48, 2).Value = "Liechtenstein"
$ws.Cells.Item(49, 2).Value = "Lithuania"
$ws.Cells.Item(50, 2).Value = "Luxembourg"
$ws.Cells.Item(51, 2).Value = "Macedonia"
$ws.Cells.Item(52, 2).Value = "Malta"
$ws.Cells.Item(53, 2).Value = "Moldova"
$ws.Cells.Item(54, 2).Value = "Monaco"
$ws.Cells.Item(55, 2).Value = "Morocco"
$ws.Cells.Item(56, 2).Value = "Netherlands"
$ws.Cells.Item(57, 2).Value = "Norway"
$ws.Cells.Item(58, 2).Value = "Oman"
$ws.Cells.Item(59, 2).Value = "Palestine"
$ws.Cells.Item(60, 2).Value = "Poland"
$ws.Cells.Item(61, 2).Value = "Portugal"
$ws.Cells.Item(62, 2).Value = "Qatar"
$ws.Cells.Item(63, 2).Value = "Romania"
$ws.Cells.Item(64, 2).Value = "Russia"
$ws.Cells.Item(65, 2).Value = "San Marino"
$ws.Cells.Item(66, 2).Value = "Saudi Arabia"
$ws.Cells.Item(67, 2).Value = "Montenegro"
$ws.Cells.Item(68, 2).Value = "Slovakia"
$ws.Cells.Item(69, 2).Value = "Slovenia"
$ws.Cells.Item(70, 2).Value = "Spain"
$ws.Cells.Item(71, 2).Value = "Sweden"
$ws.Cells.Item(72, 2).Value = "Switzerland"
$ws.Cells.Item(73, 2).Value = "Syria"
$ws.Cells.Item(74, 2).Value = "Tunisia"
$ws.Cells.Item(75, 2).Value = "Turkey"
$ws.Cells.Item(76, 2).Value = "Ukraine"
$ws.Cells.Item(77, 2).Value = "United Arab Emirates"
$ws.Cells.Item(78, 2).Value = "United Kingdom"
$ws.Cells.Item(79, 2).Value = "Yemen"
$ws.Cells.Item(80, 2).Value = "Jordan"
$ws.Cells.Item(83, 2).Value = "Turkmenistan"
$ws.Cells.Item(85, 2).Value = "Norway"
$ws.Cells.Item(86, 2).Value = "China"

# Match header B1 style (bold) to A1
$null = $ws.Range("A1").Copy()
$null = $ws.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Make Country the active/selected sheet with B1 selected (was SampleContext before)
$null = $ws.Activate()
$null = $ws.Range("B1").Select()
